$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Rajasthan Royals", "Chennai Super Kings", "Tom Curran ", "10", "9", "1", "0", "111.11"),
    @(" Dubai (DSC)", " September 30 2020", "KKR won by 37 runs", "Rajasthan Royals", "Kolkata Knight Riders", "Tom Curran ", "54", "36", "2", "3", "150.00"),
    @(" Sharjah", " September 27 2020", "Royals won by 4 wickets (with 3 balls remaining)", "Rajasthan Royals", "Kings XI Punjab", "Tom Curran ", "4", "1", "1", "0", "400.00"),
    @(" Abu Dhabi", " October 06 2020", "Mumbai won by 57 runs", "Rajasthan Royals", "Mumbai Indians", "Tom Curran ", "15", "16", "1", "0", "93.75")
)

$startRow = 6
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($c -ge 7) {
            # totalRuns/totalBalls/total4s/total6s/sr look numeric - force text storage
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$c - 1]
    }
}
